$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the age-range buckets with a coarser set of 3 ranges and update
# the corresponding Female/Male counts (rows 2-4).
$ws.Range("A2").Value = "12-15 years"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

$ws.Range("A3").Value = "15-20 years"
$ws.Range("B3").Value = 45
$ws.Range("C3").Value = 14

$ws.Range("A4").Value = "20-30 years"
$ws.Range("B4").Value = 30
$ws.Range("C4").Value = 11

# Remove the now-unused rows 5-8 (old age-range buckets) entirely so the
# sheet's used range shrinks back down to A1:C4.
$ws.Range("A5:C8").EntireRow.Delete()
